# Add data for 2022-03-10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header label to reflect the new "through" date.
$ws.Name = "Through 2022-03-02"
$ws.Range("B1").Value = "March 2022 (through March 02)"

# Add the new incident counts (value 1) for the relevant neighborhood/year cells.
$ws.Range("E5").Value = 1    # Austin
$ws.Range("N23").Value = 1   # Bridgeport
$ws.Range("E27").Value = 1   # Calumet Heights
$ws.Range("B31").Value = 1   # Near South Side
$ws.Range("H44").Value = 1   # Grand Boulevard
$ws.Range("K50").Value = 1   # Albany Park
$ws.Range("E54").Value = 1   # Avalon Park
